$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-29 Monday", "2025-12-30 Tuesday"),
    @("895×9=", "741×3="),
    @("464×7=", "299×7="),
    @("695×3=", "556×4="),
    @("228×8=", "758×4="),
    @("363×5=", "519×8="),
    @("627×9=", "633×4="),
    @("769×4=", "941×6="),
    @("598×9=", "541×3="),
    @("840×9=", "716×2="),
    @("221×4=", "834×5="),
    @("341×4=", "992×5="),
    @("713×8=", "250×9="),
    @("433×3=", "539×2="),
    @("413×3=", "513×8="),
    @("717×6=", "568×2="),
    @("412×8=", "575×8="),
    @("211×3=", "622×9="),
    @("199×4=", "500×4="),
    @("567×9=", "594×6="),
    @("942×5=", "949×8="),
    @("754×9=", "350×5="),
    @("544×8=", "937×7="),
    @("334×4=", "966×7="),
    @("847×7=", "425×5="),
    @("836×3=", "857×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
